$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (B and D).
# Excel's ColumnWidth property pads by ~0.8333 chars (the 5px/MDW grid-line
# allowance) before it is persisted as the <col width="..."> attribute, so we
# back that padding out here to land on the exact stored widths of 23 and 33.
$ws.Columns.Item(2).ColumnWidth = 23 - 0.8333333333
$ws.Columns.Item(4).ColumnWidth = 33 - 0.8333333333

# New data set (rows 2 through 21), replacing the old rows 2-7
$data = @(
    @("https://www.tiktok.com/@notellieyong/video/7544270475828940037", "notellieyong", "7544270475828940037", "Video by @notellieyong", "dance"),
    @("https://www.tiktok.com/@vlea20/video/7537390559388814614", "vlea20", "7537390559388814614", "Video by @vlea20", "dance"),
    @("https://www.tiktok.com/@el1epretty/video/7542556415735483655", "el1epretty", "7542556415735483655", "Video by @el1epretty", "dance"),
    @("https://www.tiktok.com/@77.valentinaaa/video/7541328315085524246", "77.valentinaaa", "7541328315085524246", "Video by @77.valentinaaa", "dance"),
    @("https://www.tiktok.com/@maligoshik/video/7541181053327740168", "maligoshik", "7541181053327740168", "Video by @maligoshik", "dance"),
    @("https://www.tiktok.com/@shaniaandmads/video/7542457398993145095", "shaniaandmads", "7542457398993145095", "Video by @shaniaandmads", "dance"),
    @("https://www.tiktok.com/@n_clarissa/video/7544383024348122390", "n_clarissa", "7544383024348122390", "Video by @n_clarissa", "dance"),
    @("https://www.tiktok.com/@evaforevahh/video/7544114446122224904", "evaforevahh", "7544114446122224904", "Video by @evaforevahh", "dance"),
    @("https://www.tiktok.com/@celynbrook.dance/video/7542963317246135574", "celynbrook.dance", "7542963317246135574", "Video by @celynbrook.dance", "dance"),
    @("https://www.tiktok.com/@wnyashclips2.0/video/7542195253478984982", "wnyashclips2.0", "7542195253478984982", "Video by @wnyashclips2.0", "dance"),
    @("https://www.tiktok.com/@zhurtik/video/7512543844957687045", "zhurtik", "7512543844957687045", "Video by @zhurtik", "dance"),
    @("https://www.tiktok.com/@ari.5369/video/7539280922726714646", "ari.5369", "7539280922726714646", "Video by @ari.5369", "dance"),
    @("https://www.tiktok.com/@el1epretty/video/7540703940527017223", "el1epretty", "7540703940527017223", "Video by @el1epretty", "dance"),
    @("https://www.tiktok.com/@lxttikaem1ly/video/7541740945637936407", "lxttikaem1ly", "7541740945637936407", "Video by @lxttikaem1ly", "dance"),
    @("https://www.tiktok.com/@baddiesofticktok/video/7540557654746270998", "baddiesofticktok", "7540557654746270998", "Video by @baddiesofticktok", "dance"),
    @("https://www.tiktok.com/@cocolu.xx/video/7523461069449301278", "cocolu.xx", "7523461069449301278", "Video by @cocolu.xx", "dance"),
    @("https://www.tiktok.com/@pla_neii/video/7543961522771430663", "pla_neii", "7543961522771430663", "Video by @pla_neii", "dance"),
    @("https://www.tiktok.com/@sadieemckennaa/video/7404604972358241542", "sadieemckennaa", "7404604972358241542", "Video by @sadieemckennaa", "dance"),
    @("https://www.tiktok.com/@dancingbabies_t.and.j/video/7514894090413231406", "dancingbabies_t.and.j", "7514894090413231406", "Video by @dancingbabies_t.and.j", "dance"),
    @("https://www.tiktok.com/@sara.guglielmetto/video/7540945564918107414", "sara.guglielmetto", "7540945564918107414", "Video by @sara.guglielmetto", "dance")
)

$rowCount = $data.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    # Leading apostrophe forces the numeric-looking Video ID to be stored as text
    # (matches the source data, which keeps full precision rather than a float).
    $ws.Cells.Item($r, 3).Value = "'" + $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# Reset style on the Video ID column back to Normal so no extra quote-prefix
# formatting/style reference is left behind on the cells.
$ws.Range("C2:C" + ($rowCount + 1)).Style = "Normal"
